$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "+ voor bedrijven gelden andere voorwaarden."
#    -> "+ voor bedrijven: " (run 1)  +  "10 weken, voor burgers: 12 weken" (run 2)
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("+ voor bedrijven gelden andere voorwaarden.")
if ($found1) {
    $r1.Text = "+ voor bedrijven: "
    $r1.Collapse(0)
    $r1.InsertAfter("10 weken, voor burgers: 12 weken")
}

# ---------------------------------------------------------------------------
# 2) Remove the first of the two empty paragraphs that follow
#    "3. Nadat een kwijtschelding is toegewezen kan men pas na 16 dagen
#     een regeling aanvragen."
# ---------------------------------------------------------------------------
$idx2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match [regex]::Escape("3. Nadat een kwijtschelding is toegewezen kan men pas na 16 dagen een regeling aanvragen.")) {
        $idx2 = $i
        break
    }
}
$found2 = ($idx2 -ne $null)
if ($found2) {
    $hostPara = $d.Paragraphs.Item($idx2)
    $emptyPara = $hostPara.Next()
    $emptyPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) "1. binnen 4 maanden beslissing"
#    -> "1. binnen " (run 1)  +  "3 maanden. Niet alles compleet? Dan 6 maanden" (run 2)
# ---------------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("1. binnen 4 maanden beslissing")
if ($found3) {
    $r3.Text = "1. binnen "
    $r3.Collapse(0)
    $r3.InsertAfter("3 maanden. Niet alles compleet? Dan 6 maanden")
}

# ---------------------------------------------------------------------------
# 4) Add a new paragraph "(veelgestelde vragen)" right after
#    "5. mogelijkheid van kwijtschelding 3 maanden na betaling van laatste
#     bedrag aanslag"
# ---------------------------------------------------------------------------
$idx4 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match [regex]::Escape("5. mogelijkheid van kwijtschelding 3 maanden na betaling van laatste bedrag aanslag")) {
        $idx4 = $i
        break
    }
}
$found4 = ($idx4 -ne $null)
if ($found4) {
    $hostPara4 = $d.Paragraphs.Item($idx4)
    $hostPara4.Range.InsertParagraphAfter()
    $newPara4 = $d.Paragraphs.Item($idx4 + 1)
    $newPara4.Range.Text = "(veelgestelde vragen)"
}

Write-Output "done: found1=$found1 found2=$found2 found3=$found3 found4=$found4"
